$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new submission data
$ws.Range("A2").Value = "BJ"
$ws.Range("B2").Value = "CORADO"
$ws.Range("C2").Value = "Fútbol"
$ws.Range("D2").Value = "Masculino"
$ws.Range("E2").Value = "Guatemala"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "Toyota"

# Remove the old second data row (row 3) and the trailing blank styled rows (4-6)
$ws.Range("A3:H6").Delete()
